$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update title (Volume/Number) and date-range shared strings (rich text runs) ---
$ws.Range("A8").Value = "Volume 32   Number  49"
$ws.Range("C9").Value = "Report Covering the Week  12/1/2025  Through  12/7/2025"

# --- Crime-statistics table updates (rows 15-31) ---
# Cells needing a style/type change (text <-> number): copy style from a stable donor cell, then set value
$ws.Range("J15").Copy($ws.Range("C15"))
$ws.Range("C15").Value = 2
$ws.Range("J15").Copy($ws.Range("F15"))
$ws.Range("F15").Value = 2
$ws.Range("A15").Copy($ws.Range("G15"))
$ws.Range("G15").Value = "0"
$ws.Range("A15").Copy($ws.Range("H15"))
$ws.Range("H15").Value = "***.*"
$ws.Range("J15").Copy($ws.Range("C18"))
$ws.Range("C18").Value = 2
$ws.Range("J15").Copy($ws.Range("C27"))
$ws.Range("C27").Value = 2
$ws.Range("A15").Copy($ws.Range("D27"))
$ws.Range("D27").Value = "0"
$ws.Range("A15").Copy($ws.Range("E27"))
$ws.Range("E27").Value = "***.*"
$ws.Range("A15").Copy($ws.Range("D28"))
$ws.Range("D28").Value = "0"
$ws.Range("A15").Copy($ws.Range("E28"))
$ws.Range("E28").Value = "***.*"
$ws.Range("A15").Copy($ws.Range("C29"))
$ws.Range("C29").Value = "0"
$ws.Range("J15").Copy($ws.Range("D29"))
$ws.Range("D29").Value = 2
$ws.Range("L22").Copy($ws.Range("E29"))
$ws.Range("E29").Value = -100
$ws.Range("A15").Copy($ws.Range("C30"))
$ws.Range("C30").Value = "0"
$ws.Range("J15").Copy($ws.Range("D30"))
$ws.Range("D30").Value = 1
$ws.Range("L22").Copy($ws.Range("E30"))
$ws.Range("E30").Value = -100

# Cells needing only a value update (style/type unchanged)
$ws.Range("I15").Value = 33
$ws.Range("K15").Value = -2.941176470588
$ws.Range("L15").Value = -5.714285714285
$ws.Range("M15").Value = 106.25
$ws.Range("N15").Value = -23.255813953488
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 20
$ws.Range("G16").Value = 22
$ws.Range("H16").Value = -9.090909090909
$ws.Range("I16").Value = 188
$ws.Range("J16").Value = 289
$ws.Range("K16").Value = -34.948096885813
$ws.Range("L16").Value = -27.969348659003
$ws.Range("M16").Value = -38.762214983713
$ws.Range("N16").Value = -88.308457711442
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = -16.666666666666
$ws.Range("F17").Value = 29
$ws.Range("G17").Value = 43
$ws.Range("H17").Value = -32.558139534883
$ws.Range("I17").Value = 508
$ws.Range("J17").Value = 549
$ws.Range("K17").Value = -7.468123861566
$ws.Range("L17").Value = 7.172995780590
$ws.Range("M17").Value = 92.424242424242
$ws.Range("N17").Value = -29.050279329608
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -71.428571428571
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 22
$ws.Range("H18").Value = -59.090909090909
$ws.Range("I18").Value = 210
$ws.Range("J18").Value = 195
$ws.Range("K18").Value = 7.692307692307
$ws.Range("L18").Value = 2.439024390243
$ws.Range("M18").Value = -54.048140043763
$ws.Range("N18").Value = -90.353697749196
$ws.Range("C19").Value = 20
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = 25
$ws.Range("F19").Value = 67
$ws.Range("G19").Value = 58
$ws.Range("H19").Value = 15.517241379310
$ws.Range("I19").Value = 885
$ws.Range("J19").Value = 934
$ws.Range("K19").Value = -5.246252676659
$ws.Range("L19").Value = 1.490825688073
$ws.Range("M19").Value = 43.435980551053
$ws.Range("N19").Value = -15.954415954416
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 28
$ws.Range("G20").Value = 22
$ws.Range("H20").Value = 27.272727272727
$ws.Range("I20").Value = 334
$ws.Range("J20").Value = 318
$ws.Range("K20").Value = 5.031446540880
$ws.Range("L20").Value = -18.932038834951
$ws.Range("M20").Value = 32.539682539682
$ws.Range("N20").Value = -86.834844304296
$ws.Range("C21").Value = 47
$ws.Range("D21").Value = 45
$ws.Range("E21").Value = 4.444444444444
$ws.Range("F21").Value = 155
$ws.Range("G21").Value = 168
$ws.Range("H21").Value = -7.738095238095
$ws.Range("I21").Value = 2160
$ws.Range("J21").Value = 2322
$ws.Range("K21").Value = -6.976744186046
$ws.Range("L21").Value = -4.719894133215
$ws.Range("M21").Value = 12.382934443288
$ws.Range("N21").Value = -73.529411764705
$ws.Range("J22").Value = 44
$ws.Range("K22").Value = -34.090909090909
$ws.Range("M22").Value = -6.451612903225
$ws.Range("C23").Value = 4
$ws.Range("D23").Value = 6
$ws.Range("F23").Value = 14
$ws.Range("G23").Value = 16
$ws.Range("H23").Value = -12.5
$ws.Range("I23").Value = 202
$ws.Range("J23").Value = 208
$ws.Range("K23").Value = -2.884615384615
$ws.Range("L23").Value = -12.931034482758
$ws.Range("M23").Value = 33.774834437086
$ws.Range("C24").Value = 47
$ws.Range("D24").Value = 40
$ws.Range("E24").Value = 17.5
$ws.Range("G24").Value = 206
$ws.Range("H24").Value = -0.970873786407
$ws.Range("I24").Value = 2177
$ws.Range("J24").Value = 2199
$ws.Range("K24").Value = -1.000454752160
$ws.Range("L24").Value = -2.332884701659
$ws.Range("M24").Value = 55.5
$ws.Range("C25").Value = 22
$ws.Range("D25").Value = 20
$ws.Range("E25").Value = 10
$ws.Range("F25").Value = 92
$ws.Range("G25").Value = 100
$ws.Range("H25").Value = -8
$ws.Range("I25").Value = 1053
$ws.Range("J25").Value = 1139
$ws.Range("K25").Value = -7.550482879719
$ws.Range("L25").Value = 3.134182174338
$ws.Range("C26").Value = 16
$ws.Range("D26").Value = 16
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 58
$ws.Range("G26").Value = 54
$ws.Range("H26").Value = 7.407407407407
$ws.Range("I26").Value = 813
$ws.Range("J26").Value = 868
$ws.Range("K26").Value = -6.336405529953
$ws.Range("L26").Value = -0.853658536585
$ws.Range("M26").Value = -2.166064981949
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 200
$ws.Range("I27").Value = 41
$ws.Range("K27").Value = -16.326530612244
$ws.Range("L27").Value = -26.785714285714
$ws.Range("C28").Value = 2
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 86
$ws.Range("K28").Value = 3.614457831325
$ws.Range("L28").Value = -21.100917431192
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = -50
$ws.Range("J29").Value = 9
$ws.Range("K29").Value = 33.333333333333
$ws.Range("F30").Value = 1
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 7
$ws.Range("K30").Value = 57.142857142857
$ws.Range("L31").Value = -73.684210526315
